# Auto-generated script applying numeric updates per the diff
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Table S1 - Plasticity AIC")
$ws1.Range("H2").Value = 0.339921366849096
$ws1.Range("I2").Value = 0.824553438127845
$ws1.Range("J2").Value = 0.29777176529434
$ws1.Range("H3").Value = 0.270712906666334
$ws1.Range("I3").Value = 0.880517357628651
$ws1.Range("J3").Value = 0.309990520594351
$ws1.Range("F4").Value = 0.511
$ws1.Range("H4").Value = 0.272155096268899
$ws1.Range("I4").Value = 0.874876231377499
$ws1.Range("J4").Value = 0.307314521780736
$ws1.Range("H5").Value = 0.282067343982829
$ws1.Range("I5").Value = 0.855828398264632
$ws1.Range("J5").Value = 0.294010052236142
$ws1.Range("H6").Value = 0.252287909332907
$ws1.Range("I6").Value = 0.933847193844437
$ws1.Range("J6").Value = 0.322867004307088
$ws1.Range("H7").Value = 0.308813845572371
$ws1.Range("I7").Value = 0.936882341333349
$ws1.Range("J7").Value = 0.331059051431514
$ws1.Range("H8").Value = 0.252780623700835
$ws1.Range("I8").Value = 0.931849378942401
$ws1.Range("J8").Value = 0.323013483263546
$ws1.Range("D9").Value = 110.7
$ws1.Range("E9").Value = 128.3
$ws1.Range("F9").Value = 0.427
$ws1.Range("G9").Value = 0.347
$ws1.Range("H9").Value = 0.12264193082729
$ws1.Range("I9").Value = 1.00977970615862
$ws1.Range("J9").Value = 0.335976868226673
$ws1.Range("D10").Value = 106
$ws1.Range("E10").Value = 118.6
$ws1.Range("F10").Value = 0.341
$ws1.Range("G10").Value = 0.292
$ws1.Range("H10").Value = 0.069512358021633
$ws1.Range("I10").Value = 1.12393872070358
$ws1.Range("J10").Value = 0.353702412571083
$ws1.Range("K10").Value = 44.8
$ws1.Range("D11").Value = 106.9
$ws1.Range("E11").Value = 119.5
$ws1.Range("F11").Value = 0.313
$ws1.Range("G11").Value = 0.271
$ws1.Range("H11").Value = 0.0572627646077092
$ws1.Range("I11").Value = 1.11383231952833
$ws1.Range("J11").Value = 0.352824835342725
$ws1.Range("K11").Value = 36.3
$ws1.Range("D12").Value = 102.8
$ws1.Range("E12").Value = 111.6
$ws1.Range("F12").Value = 0.268
$ws1.Range("G12").Value = 0.224
$ws1.Range("H12").Value = 0.056048071109652
$ws1.Range("I12").Value = 1.19677569320301
$ws1.Range("J12").Value = 0.36408716814135
$ws1.Range("K12").Value = 28.6
$ws1.Range("H13").Value = 0.405142466313655
$ws1.Range("I13").Value = 0.710112914902594
$ws1.Range("J13").Value = 0.231080474175141
$ws1.Range("H14").Value = 0.405616901269354
$ws1.Range("I14").Value = 0.709947971925215
$ws1.Range("J14").Value = 0.230951400774421
$ws1.Range("H15").Value = 0.409641426216089
$ws1.Range("I15").Value = 0.70697164728988
$ws1.Range("J15").Value = 0.229916347254973
$ws1.Range("H16").Value = 0.40760856417087
$ws1.Range("I16").Value = 0.727645208726944
$ws1.Range("J16").Value = 0.236972613071863
$ws1.Range("H17").Value = 0.393286214699538
$ws1.Range("I17").Value = 0.73324439340205
$ws1.Range("J17").Value = 0.237295448480225
$ws1.Range("H18").Value = 0.396708536396551
$ws1.Range("I18").Value = 0.732981216183524
$ws1.Range("J18").Value = 0.238221370362628
$ws1.Range("H19").Value = 0.394017895328747
$ws1.Range("I19").Value = 0.734044282144376
$ws1.Range("J19").Value = 0.23735001040535

$ws2 = $wb.Worksheets.Item("Table S2 - PERMANOVA")
$ws2.Range("F3").Value = 0.05929
$ws2.Range("F9").Value = 0.22785
$ws2.Range("F12").Value = 0.46036
$ws2.Range("F13").Value = 0.00133

$ws3 = $wb.Worksheets.Item("Table S3 - Plasticity GLM")
$ws3.Range("D2").Value = 187.01
$ws3.Range("D3").Value = 1.72
$ws3.Range("E3").Value = 0.085
$ws3.Range("D4").Value = 66.95
$ws3.Range("D5").Value = 39.32
$ws3.Range("D6").Value = 78.75
$ws3.Range("D7").Value = 0.31
$ws3.Range("E7").Value = 0.753
$ws3.Range("D8").Value = -136.81
$ws3.Range("D9").Value = -79.33
$ws3.Range("D10").Value = -58.48
$ws3.Range("B19").Value = 1.033
$ws3.Range("C19").Value = 0.125
$ws3.Range("D19").Value = 8.25
$ws3.Range("B20").Value = -0.045
$ws3.Range("C20").Value = 0.115
$ws3.Range("D20").Value = -0.39
$ws3.Range("E20").Value = 0.694
$ws3.Range("B21").Value = 0.033
$ws3.Range("C21").Value = 0.079
$ws3.Range("D21").Value = 0.42
$ws3.Range("E21").Value = 0.676
$ws3.Range("B22").Value = 0.124
$ws3.Range("C22").Value = 0.082
$ws3.Range("D22").Value = 1.51
$ws3.Range("E22").Value = 0.131
$ws3.Range("B23").Value = 0.262
$ws3.Range("C23").Value = 0.069
$ws3.Range("D23").Value = 3.82
$ws3.Range("B24").Value = 0.493
$ws3.Range("B25").Value = 0.145

$ws4 = $wb.Worksheets.Item("Table S4 - Species PERMANOVA")
$ws4.Range("F3").Value = 0.0986
$ws4.Range("F4").Value = 0.004
$ws4.Range("F7").Value = 0.01865
$ws4.Range("F8").Value = 0.00533

$ws5 = $wb.Worksheets.Item("Table S5 - HostVsymb PERMANOVA")
$ws5.Range("F3").Value = 0.00666
$ws5.Range("K3").Value = 0.0966
$ws5.Range("F4").Value = 0.56762
$ws5.Range("K4").Value = 0.00266
$ws5.Range("F7").Value = 0.27515
$ws5.Range("K7").Value = 0.28115
$ws5.Range("F9").Value = 0.16855
$ws5.Range("K9").Value = 0.00133
$ws5.Range("F12").Value = 0.01332
$ws5.Range("K12").Value = 0.00067
$ws5.Range("F13").Value = 0.09127
$ws5.Range("F14").Value = 0.20253
$ws5.Range("K14").Value = 0.46636

$ws6 = $wb.Worksheets.Item("Table S6 - HostVsymb Plast AIC")
$ws6.Range("G2").Value = 0.163808918846475
$ws6.Range("H2").Value = 0.884761957169536
$ws6.Range("I2").Value = 0.445413568051623
$ws6.Range("G3").Value = 0.155601596236995
$ws6.Range("H3").Value = 0.913117790532039
$ws6.Range("I3").Value = 0.45931696402324
$ws6.Range("G4").Value = 0.155601596236034
$ws6.Range("H4").Value = 0.913117790532143
$ws6.Range("I4").Value = 0.459316964023365
$ws6.Range("G5").Value = 0.155601596236034
$ws6.Range("H5").Value = 0.913117790532143
$ws6.Range("I5").Value = 0.459316964023365
$ws6.Range("F6").Value = 0.212
$ws6.Range("G6").Value = 0.140275722020262
$ws6.Range("H6").Value = 0.935551604873005
$ws6.Range("I6").Value = 0.476818047697685
$ws6.Range("G7").Value = 0.163720683028945
$ws6.Range("H7").Value = 0.950425796937665
$ws6.Range("I7").Value = 0.488694314791684
$ws6.Range("J7").Value = 59.1
$ws6.Range("G8").Value = 0.161404588856711
$ws6.Range("H8").Value = 0.950207028246836
$ws6.Range("I8").Value = 0.488554386709045
$ws6.Range("G9").Value = 0.144860307744428
$ws6.Range("H9").Value = 0.948916968928776
$ws6.Range("I9").Value = 0.477395026003451
$ws6.Range("G10").Value = 0.157143513626403
$ws6.Range("H10").Value = 0.982507593585701
$ws6.Range("I10").Value = 0.505803686093402
$ws6.Range("G11").Value = 0.121607431312658
$ws6.Range("H11").Value = 1.01492820585159
$ws6.Range("I11").Value = 0.513921291826242
$ws6.Range("G12").Value = 0.122648785796666
$ws6.Range("H12").Value = 1.00853413602528
$ws6.Range("I12").Value = 0.508488836608921

$ws7 = $wb.Worksheets.Item("Table S7 - HostVsymb Plast GLM")
$ws7.Range("B3").Value = 0.199
$ws7.Range("E3").Value = 0.492
$ws7.Range("D4").Value = -1.46
$ws7.Range("B6").Value = 0.234
$ws7.Range("E6").Value = 0.273
$ws7.Range("B7").Value = 0.005
$ws7.Range("E7").Value = 0.974
$ws7.Range("E8").Value = 0.895
$ws7.Range("D9").Value = -1.96
$ws7.Range("E9").Value = 0.05
$ws7.Range("B10").Value = 0.223
$ws7.Range("E10").Value = 0.47
$ws7.Range("D11").Value = 2.53
$ws7.Range("B13").Value = -0.568
$ws7.Range("E13").Value = 0.094
$ws7.Range("B14").Value = -0.45
$ws7.Range("B15").Value = 0.051
$ws7.Range("E15").Value = 0.833
$ws7.Range("D16").Value = -0.53
$ws7.Range("E16").Value = 0.593
$ws7.Range("E17").Value = 0.335
$ws7.Range("B18").Value = 0.072
$ws7.Range("E18").Value = 0.774
$ws7.Range("E21").Value = 0.555
$ws7.Range("E24").Value = 0.659
$ws7.Range("E25").Value = 0.154
$ws7.Range("B26").Value = -0.179
$ws7.Range("D26").Value = -0.43
$ws7.Range("E26").Value = 0.669
$ws7.Range("B27").Value = -0.18
$ws7.Range("B29").Value = -0.89
$ws7.Range("B30").Value = 0.189
$ws7.Range("E30").Value = 0.589
$ws7.Range("B31").Value = -0.607

Write-Output "Applied 195 cell updates across 7 sheets"
